# Updated cryptos list on Sat Apr 29 18:25:04 UTC 2023 with GitHub Actions
#
# Refresh the Price (D) / Volume(1h) (E) columns with the latest scrape, and
# swap a handful of adjacent rows whose relative ranking flipped (Polkadot vs
# WrappedEther, Aptos vs TheSandbox, RenderToken vs WEMIXToken, EnergySwap vs
# Cronos). A leading apostrophe is used on numeric-looking Price strings so
# Excel keeps them as text (preserving things like trailing zeros and the
# "thousands-dot" look of "29.416.75") instead of auto-coercing to a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "'29.416.75"
$ws.Range("E2").Value = "  +0.67%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "'1.913.14"

# Row 4 - TetherUSD
$ws.Range("D4").Value = "'1.010"
$ws.Range("E4").Value = "  +0.82%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'325.35"
$ws.Range("E5").Value = "  +0.75%  "

# Row 6 - USDC
$ws.Range("D6").Value = "'1.008"
$ws.Range("E6").Value = "  +0.63%  "

# Row 7 - XRP
$ws.Range("D7").Value = "'0.4805"
$ws.Range("E7").Value = "  +1.28%  "

# Row 8 - Cardano
$ws.Range("D8").Value = "'0.4068"
$ws.Range("E8").Value = "  +0.78%  "

# Row 9 - Dogecoin
$ws.Range("D9").Value = "'0.08250"
$ws.Range("E9").Value = "  +2.82%  "

# Row 10 - Polygon
$ws.Range("E10").Value = "  +2.34%  "

# Row 11 - Solana
$ws.Range("D11").Value = "'23.46"
$ws.Range("E11").Value = "  +0.65%  "

# Row 12 - was WrappedEther, now Polkadot (rows 12/13 swapped)
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").Value = "'6.034"
$ws.Range("E12").Value = "  +1.94%  "

# Row 13 - was Polkadot, now WrappedEther
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "'1.878.04"
$ws.Range("E13").Value = "  -1.03%  "

# Row 14 - Chainlink
$ws.Range("D14").Value = "'7.207"
$ws.Range("E14").Value = "  +2.60%  "

# Row 15 - Litecoin
$ws.Range("D15").Value = "'91.10"
$ws.Range("E15").Value = "  +2.35%  "

# Row 16 - TRON
$ws.Range("D16").Value = "'0.06808"
$ws.Range("E16").Value = "  +2.63%  "

# Row 17 - BinanceUSD
$ws.Range("E17").Value = "  +0.76%  "

# Row 18 - ShibaInu
$ws.Range("D18").Value = "'0.00001040"
$ws.Range("E18").Value = "  +1.17%  "

# Row 19 - Avalanche
$ws.Range("D19").Value = "'17.70"
$ws.Range("E19").Value = "  +1.21%  "

# Row 20 - Dai
$ws.Range("D20").Value = "'1.008"
$ws.Range("E20").Value = "  +0.72%  "

# Row 21 - WrappedBTC
$ws.Range("D21").Value = "'29.447.43"
$ws.Range("E21").Value = "  +0.69%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "'5.623"
$ws.Range("E22").Value = "  +2.15%  "

# Row 23 - Cosmos
$ws.Range("D23").Value = "'11.79"
$ws.Range("E23").Value = "  +0.41%  "

# Row 24 - Toncoin
$ws.Range("D24").Value = "'2.196"
$ws.Range("E24").Value = "  +1.97%  "

# Row 25 - WrappedliquidstakedEther2.0
$ws.Range("D25").Value = "'2.135.53"
$ws.Range("E25").Value = "  +1.53%  "

# Row 26 - InternetComputer(DFINITY)
$ws.Range("D26").Value = "'6.589"
$ws.Range("E26").Value = "  +11.50%  "

# Row 27 - Monero
$ws.Range("D27").Value = "'156.77"
$ws.Range("E27").Value = "  +1.46%  "

# Row 28 - EthereumClassic
$ws.Range("D28").Value = "'20.06"
$ws.Range("E28").Value = "  +1.83%  "

# Row 29 - LidoDAOToken
$ws.Range("D29").Value = "'2.105"
$ws.Range("E29").Value = "  +1.23%  "

# Row 30 - BitcoinCash
$ws.Range("D30").Value = "'120.14"
$ws.Range("E30").Value = "  +2.14%  "

# Row 31 - ImmutableX
$ws.Range("D31").Value = "'1.020"
$ws.Range("E31").Value = "  -0.70%  "

# Row 32 - Stellar
$ws.Range("D32").Value = "'0.09558"
$ws.Range("E32").Value = "  +1.65%  "

# Row 33 - Filecoin
$ws.Range("D33").Value = "'5.569"
$ws.Range("E33").Value = "  +4.40%  "

# Row 34 - HuobiToken
$ws.Range("D34").Value = "'3.553"
$ws.Range("E34").Value = "  +0.72%  "

# Row 35 - ARBITRUM
$ws.Range("D35").Value = "'1.366"
$ws.Range("E35").Value = "  -0.28%  "

# Row 36 - VeChain
$ws.Range("E36").Value = "  +1.56%  "

# Row 37 - Hedera
$ws.Range("D37").Value = "'0.06119"
$ws.Range("E37").Value = "  +1.50%  "

# Row 38 - TrustWalletToken
$ws.Range("D38").Value = "'1.183"
$ws.Range("E38").Value = "  +1.23%  "

# Row 39 - FraxShare
$ws.Range("D39").Value = "'8.043"
$ws.Range("E39").Value = "  +1.74%  "

# Row 40 - was Aptos, now TheSandbox (rows 40/41 swapped)
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").Value = "'0.5965"
$ws.Range("E40").Value = "  +2.47%  "

# Row 41 - was TheSandbox, now Aptos
$ws.Range("B41").Value = "Aptos"
$ws.Range("C41").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D41").Value = "'10.83"
$ws.Range("E41").Value = "  +7.91%  "

# Row 42 - Algorand
$ws.Range("E42").Value = "  +1.14%  "

# Row 43 - was RenderToken, now WEMIXToken (rows 43/44 swapped)
$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").Value = "'1.281"
$ws.Range("E43").Value = "  -0.63%  "

# Row 44 - was WEMIXToken, now RenderToken
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").Value = "'2.394"
$ws.Range("E44").Value = "  +1.27%  "

# Row 45 - was EnergySwap, now Cronos (rows 45/46 swapped)
$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D45").Value = "'0.07608"
$ws.Range("E45").Value = "  -1.29%  "

# Row 46 - was Cronos, now EnergySwap
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'12.42"
$ws.Range("E46").Value = "  +1.55%  "

# Row 47 - Decentraland
$ws.Range("D47").Value = "'0.5569"
$ws.Range("E47").Value = "  +1.62%  "

# Row 48 - NEARProtocol
$ws.Range("D48").Value = "'1.952"
$ws.Range("E48").Value = "  +2.33%  "

# Row 49 - Quant
$ws.Range("D49").Value = "'117.41"
$ws.Range("E49").Value = "  +3.98%  "

# Row 50 - MXToken
$ws.Range("D50").Value = "'2.429"
$ws.Range("E50").Value = "  +4.24%  "

# Row 51 - Aave
$ws.Range("D51").Value = "'72.22"
$ws.Range("E51").Value = "  +1.27%  "
